# "Data source corrected and updated"
#
# The J and K columns held placeholder/mismatched values (including two
# shared-string labels "r"/"s" in row 1 instead of numbers). Re-point them
# at the corrected data source: the whole J column becomes a constant 1 and
# the whole K column becomes a constant 0.6 for every data row (1-51).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J1:J51").Value = 1
$ws.Range("K1:K51").Value = 0.6

# Leave the selection on the refreshed column, matching the author's
# post-edit selection state (K1 active, K1:K51 selected).
$ws.Range("K1:K51").Select()
